$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2 through 28 from 45224 (2023-10-25)
# to 45233 (2023-11-03), keeping existing number formatting intact.
$ws.Range("C2:C28").Value = 45233
